$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.480.02'
$ws.Range("E2").Value = '  -3.03%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.775.57'
$ws.Range("E3").Value = '  -2.01%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.31%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.31'
$ws.Range("E5").Value = '  -0.81%  '

$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9998'
$ws.Range("E6").Value = '  -0.27%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4284'
$ws.Range("E7").Value = '  +1.64%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3641'
$ws.Range("E8").Value = '  +2.16%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07155'
$ws.Range("E9").Value = '  -0.04%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8479'
$ws.Range("E10").Value = '  +0.03%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.50'
$ws.Range("E11").Value = '  +1.48%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.770.12'
$ws.Range("E12").Value = '  -3.09%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.454'
$ws.Range("E13").Value = '  +1.11%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.261'
$ws.Range("E14").Value = '  -1.47%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06868'
$ws.Range("E15").Value = '  -0.80%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.004'
$ws.Range("E16").Value = '  -0.23%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '78.98'
$ws.Range("E17").Value = '  -2.95%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008646'
$ws.Range("E18").Value = '  -2.06%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9982'
$ws.Range("E19").Value = '  -0.40%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.97'
$ws.Range("E20").Value = '  -1.36%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.485.34'
$ws.Range("E21").Value = '  -4.88%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.117'
$ws.Range("E22").Value = '  +0.37%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.10'
$ws.Range("E23").Value = '  +1.47%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.989.77'
$ws.Range("E24").Value = '  -3.75%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.32'
$ws.Range("E25").Value = '  -0.79%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.871'
$ws.Range("E26").Value = '  -4.84%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.02'
$ws.Range("E27").Value = '  -1.26%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.082'
$ws.Range("E28").Value = '  -0.36%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.82'
$ws.Range("E29").Value = '  +0.34%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.810'
$ws.Range("E30").Value = '  +4.50%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08930'
$ws.Range("E31").Value = '  +0.25%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7297'
$ws.Range("E32").Value = '  -1.83%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.130'
$ws.Range("E33").Value = '  +1.66%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.330'
$ws.Range("E34").Value = '  -3.48%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.743'
$ws.Range("E35").Value = '  -6.56%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9990'
$ws.Range("E36").Value = '  -0.35%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.107'
$ws.Range("E37").Value = '  +3.10%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05159'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01891'
$ws.Range("E39").Value = '  -0.87%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4933'
$ws.Range("E40").Value = '  -1.27%  '

$ws.Range("E41").Value = '  -1.77%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.630'
$ws.Range("E42").Value = '  -5.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.326'
$ws.Range("E43").Value = '  +0.29%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.031'
$ws.Range("E44").Value = '  -2.33%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '105.01'
$ws.Range("E45").Value = '  -0.24%  '

$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9991'
$ws.Range("E46").Value = '  -0.27%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.13'
$ws.Range("E47").Value = '  -2.06%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.642'
$ws.Range("E48").Value = '  +2.33%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06205'
$ws.Range("E49").Value = '  -3.25%  '

$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4481'
$ws.Range("E50").Value = '  -2.64%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.731'
$ws.Range("E51").Value = '  +2.69%  '
